$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Preserve B2's current number format (quotePrefix xf) across the value
# change by round-tripping it through a scratch cell far off the used range.
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# Row 1 - headers
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "AdminUsername"
$ws.Range("B1").Value = "AdminPassword"
$ws.Range("C1").Value = "StandardUser"
$ws.Range("D1").Value = "StandardPassword"
$ws.Range("E1").Value = "NonRegisterUser"
$ws.Range("F1").Value = "InvalidPassword"

# New header cells need the same shaded header format as the existing ones.
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# Row 2 - values
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "admin@cgi.com"
$ws.Range("B2").Value = "adminpw"
$ws.Range("C2").Value = "sravan.neppalli@cgi.com"
$ws.Range("D2").Value = "ddff444"
$ws.Range("E2").Value = "venkat@cgi.com"
$ws.Range("F2").Value = "venkat"

# Restore B2's original (quotePrefix) format, then drop the scratch cell.
$ws.Range("Z1").Copy()
$ws.Range("B2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Z1").Clear()

# ---------------------------------------------------------------------------
# Hyperlinks - drop the old set and recreate rId1/rId2/rId3 in column order
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:admin@cgi.com")
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:sravan.neppalli@cgi.com")
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:venkat@cgi.com")

# Hyperlinks.Add() stamps a freshly-minted xf; point these cells back at the
# shared built-in "Hyperlink" cell style instead.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("E2").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Column widths for the two new columns + width tweaks on existing ones.
# (Values are pre-compensated for this host's char<->pixel rounding so the
# stored <col width=.../> lands as close as possible to the target widths.)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 13.166666666666666
$ws.Columns.Item(3).ColumnWidth = 20.833333333333332
$ws.Columns.Item(4).ColumnWidth = 20.833333333333332
$ws.Columns.Item(5).ColumnWidth = 15.666666666666666
$ws.Columns.Item(6).ColumnWidth = 13.333333333333334

# ---------------------------------------------------------------------------
# Selection matches the new last-used cell
# ---------------------------------------------------------------------------
$ws.Range("F2").Select() | Out-Null
